$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (A1:K1) - inserted new "SOC_start" column at F, renamed former
# SOC_start (old H) to "SOC_end" (now I), shifting v/area_start/SOC_warn/
# P_charge_lambda one column to the right.
$headers = @("C_max", "P", "T_start", "T_end", "P_charge", "SOC_start", "v", "area_start", "SOC_end", "SOC_warn", "P_charge_lambda")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Data rows (A2:K21) - full refreshed dispatch/SOC dataset
$data = @(
    @(106, 3, 18, 79, 5, 0.6, 2.7, 4, 0.9, 0.4, 0.9),
    @(106, 3, 22, 64, 5, 0.6, 2.8, 2, 0.8, 0.4, 0.9),
    @(109, 4, 20, 66, 5, 0.6, 3.4, 4, 0.8, 0.4, 0.9),
    @(101, 3, 24, 70, 5, 0.6, 2.8, 3, 0.9, 0.4, 0.9),
    @(112, 3, 25, 79, 5, 0.6, 2.9, 10, 0.9, 0.4, 0.9),
    @(106, 4, 30, 74, 5, 0.5, 3.3, 4, 0.8, 0.4, 0.9),
    @(106, 3, 29, 86, 5, 0.6, 3, 3, 0.9, 0.4, 0.9),
    @(108, 3, 24, 72, 5, 0.6, 3.3, 8, 0.8, 0.4, 0.9),
    @(113, 3, 24, 77, 5, 0.7, 3.2, 4, 0.9, 0.4, 0.9),
    @(101, 3, 28, 74, 5, 0.7, 3.2, 10, 0.9, 0.4, 0.9),
    @(104, 3, 23, 73, 5, 0.6, 2.8, 12, 0.9, 0.4, 0.9),
    @(103, 4, 32, 65, 5, 0.5, 3, 1, 0.8, 0.4, 0.9),
    @(106, 3, 21, 72, 5, 0.7, 3.2, 1, 0.8, 0.4, 0.9),
    @(105, 2, 23, 71, 5, 0.6, 2.5, 9, 0.9, 0.4, 0.9),
    @(104, 3, 28, 61, 5, 0.7, 3.2, 5, 0.9, 0.4, 0.9),
    @(106, 3, 29, 72, 4, 0.6, 2.5, 6, 0.9, 0.4, 0.9),
    @(108, 3, 34, 56, 5, 0.6, 3, 2, 0.8, 0.4, 0.9),
    @(102, 3, 41, 71, 5, 0.6, 2.8, 15, 0.9, 0.4, 0.9),
    @(108, 3, 34, 66, 5, 0.7, 3.1, 9, 0.9, 0.4, 0.9),
    @(107, 2, 24, 73, 5, 0.6, 2.6, 7, 0.9, 0.4, 0.9)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
